# Add two new i18n string rows (app.urlcleaner.confirmBtn / app.urlcleaner.ruleTitle)
# to the bottom of the translation table on Sheet1, following the exact layout of
# the preceding rows (A=key, B=zh-CN, C=en-US, D/E left blank for zh-HK/fr-FR).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 54: app.urlcleaner.confirmBtn ---------------------------------
# Clone formatting (fill/border/alignment + row height) from the last existing
# data row (53) so the new row matches the table's styling exactly.
$ws.Range("A53:E53").Copy()
$ws.Range("A54:E54").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows("54").RowHeight = 20.1

$ws.Range("A54").Value = "app.urlcleaner.confirmBtn"
$ws.Range("B54").Value = "净化"
$ws.Range("C54").Value = "Clean URL"

# --- Row 55: app.urlcleaner.ruleTitle -----------------------------------
$ws.Range("A53:E53").Copy()
$ws.Range("A55:E55").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows("55").RowHeight = 20.1

$ws.Range("A55").Value = "app.urlcleaner.ruleTitle"
$ws.Range("B55").Value = "规则"
$ws.Range("C55").Value = "Rules"
